# National_Measles_Cases_Weekly_2025.xlsx - "update data and optimization function"
#
# 1) Sheet "National_Measles_Cases_Weekly": drop the imported_cases / local_cases
#    columns (C, D) -- the sheet is now just EpiWeek + total_cases -- and append
#    four new weekly rows (EpiWeek 39-42).
# 2) Sheet "Sources": the CDC NNDSS row/hyperlink is removed, the remaining
#    row is relabeled "Total" / JHU tracker url, and the sheet becomes the
#    active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # National_Measles_Cases_Weekly
$ws2 = $wb.Worksheets.Item(2)   # Sources

# --- Sheet1: remove the imported_cases (C) and local_cases (D) columns ---
$ws1.Columns.Item(3).Delete()
$ws1.Columns.Item(3).Delete()

# --- Sheet1: append the new weekly totals (EpiWeek 39-42) ---
$ws1.Range("A40").Value = 39
$ws1.Range("B40").Value = 31
$ws1.Range("A41").Value = 40
$ws1.Range("B41").Value = 15
$ws1.Range("A42").Value = 41
$ws1.Range("B42").Value = 29
$ws1.Range("A43").Value = 42
$ws1.Range("B43").Value = 18

# --- Sheet1: move the live selection down to the new last row ---
[void]$ws1.Range("E43").Select()

# --- Sheet2: drop the CDC source row/hyperlink, relabel the remaining one ---
$ws2.Hyperlinks.Delete()
$ws2.Range("A2").ClearContents()
$ws2.Range("B2").Value = ""
$ws2.Range("A1").Value = "Total"
$ws2.Range("B1").Value = "https://publichealth.jhu.edu/ivac/resources/us-measles-tracker"

# --- Sources becomes the active/selected tab ---
[void]$ws2.Activate()
